$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly entered values (cells previously empty, now set to 5)
$ws.Range("C9").Value = 5
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = 5
$ws.Range("D15").Value = 5
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 5

# Update the active selection to match the saved view state (C9)
$ws.Range("C9").Select()
